# MHD2-259: Report template and related changes for reporting on 136 genes
#
# The document contains a single outer table (one row, one cell) that wraps
# the whole clinical-context body. Two fill colours change:
#   1. The table-level shading (w:tblPr/w:shd)  FFF2CC -> ECEAF2
#   2. The heading cell's shading (w:tcPr/w:shd) E8E7EC -> ECEAF2
#
# Helper: convert a "RRGGBB" hex string into the BGR-packed OLE color value
# that Word's Shading.BackgroundPatternColor expects.
function HexToOleColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$d = $word.ActiveDocument
$newColor = HexToOleColor "ECEAF2"

$table = $d.Tables.Item(1)

# 1) Table-wide shading, stored on <w:tblPr><w:shd .../> - reached through
#    Table.Rows.Shading (applies the shading at the table-properties level).
$table.Rows.Shading.BackgroundPatternColor = $newColor

# 2) The first (heading) cell's own shading, stored on <w:tcPr><w:shd .../>.
$cell = $table.Cell(1, 1)
$cell.Shading.BackgroundPatternColor = $newColor

Write-Host "Table shading set to:" $table.Rows.Shading.BackgroundPatternColor
Write-Host "Cell shading set to:" $cell.Shading.BackgroundPatternColor
